$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0785058923029907
$ws.Range("C2").Value = 0.07879932249964856
$ws.Range("B3").Value = 25.8507665079589
$ws.Range("C3").Value = 25.84962391117613
$ws.Range("B4").Value = 151.7000171572285
$ws.Range("C4").Value = 151.6998888889504
$ws.Range("B5").Value = 0.1860813878608437
$ws.Range("C5").Value = 0.1863763212861354
$ws.Range("B6").Value = 2.2599826594633
$ws.Range("C6").Value = 2.260010699463755
$ws.Range("B7").Value = 3.446411539251389
$ws.Range("C7").Value = 3.460343982477563
$ws.Range("B8").Value = 28.34446575857872
$ws.Range("C8").Value = 28.35036153016067
$ws.Range("B9").Value = 0.9985560960564827
$ws.Range("C9").Value = 1.001990562274371
$ws.Range("B10").Value = 1.037626909440452
$ws.Range("C10").Value = 1.02993931306435
$ws.Range("B11").Value = 2.147901651422697
$ws.Range("C11").Value = 2.147528876537519
$ws.Range("B12").Value = 0.9142274572395582
$ws.Range("C12").Value = 0.9137685644164281
$ws.Range("B13").Value = 0.07819707116398543
$ws.Range("C13").Value = 0.07852662361418286
$ws.Range("B16").Value = 0.5611650814235752
$ws.Range("C16").Value = 0.5596971337829894
$ws.Range("B17").Value = 0.1760028555701231
$ws.Range("C17").Value = 0.1782905614990401
$ws.Range("B18").Value = 0.1238448749993802
$ws.Range("C18").Value = 0.1296115057471527
$ws.Range("B19").Value = 5.106580426202616
$ws.Range("C19").Value = 5.109398435461945
$ws.Range("B20").Value = -0.001518639068848169
$ws.Range("C20").Value = 0.000425050072496377
$ws.Range("B21").Value = 0.466748659620552
$ws.Range("C21").Value = 0.4728653484395555
$ws.Range("B22").Value = 28.35098618505414
$ws.Range("C22").Value = 28.35036122244466
$ws.Range("B23").Value = 29.29968616819997
$ws.Range("C23").Value = 29.29961077147871
$ws.Range("B24").Value = 0.08170680757140904
$ws.Range("C24").Value = 0.08153342587004436
$ws.Range("B25").Value = 0.4545569202934908
$ws.Range("C25").Value = 0.4549170933143945
$ws.Range("B26").Value = 1.706368158243851
$ws.Range("C26").Value = 1.707962202343462
$ws.Range("B27").Value = 3.593914212738785
$ws.Range("C27").Value = 3.593589538502072
$ws.Range("B28").Value = 12.06594184675126
$ws.Range("C28").Value = 12.06656282398755
$ws.Range("B29").Value = 33.53453594755665
$ws.Range("C29").Value = 34.14884386934217
$ws.Range("B30").Value = 73731.47787848058
$ws.Range("C30").Value = 73834.8327812543
$ws.Range("B31").Value = 6.091030113757384
$ws.Range("C31").Value = 6.093138382163727
$ws.Range("B32").Value = 87.8510884073508
$ws.Range("C32").Value = 88.07489705110397
$ws.Range("B33").Value = -0.004194292748618777
$ws.Range("C33").Value = 0.00122776659692803
